$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with newer daily COVID numbers (Date, Daily Cases,
# Confirmed, Deaths, Recovered, Discarded, Analyze) by copying the
# formatting of the last existing data row (100) down to the new rows,
# then filling in the values - mirrors "drag-fill" / copy-down editing.
$ws.Range("A100:G100").Copy()
$ws.Range("A101:G104").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$data = @(
    @(44022, 18, 1410, 74, 1055, 855, 5),
    @(44023, 11, 1424, 74, 1073, 862, 2),
    @(44024, 11, 1435, 74, 1085, 870, 4),
    @(44025, 10, 1445, 74, 1100, 877, 4)
)

$row = 101
foreach ($r in $data) {
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($row, $c).Value = $r[$c - 1]
    }
    $row = $row + 1
}

# Row 105 only carries the date-formatted (empty) cell in column A, same
# as the pattern at the bottom of the sheet before this edit.
$ws.Range("A104").Copy()
$ws.Range("A105").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("C107").Select()
